{"js": "// Update the date title paragraph from \"2025-03-13 Thursday\" to \"2025-03-14 Friday\".\nconst titleOld = \"2025-03-13 Thursday\";\nconst titleNew = \"2025-03-14 Friday\";\nconst titleResults = context.document.body.search(titleOld, { matchCase: true });\ntitleResults.load(\"text\");\nawait context.sync();\nif (titleResults.items.length > 0) {\n  titleResults.items[0].insertText(titleNew, Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// Update every math-problem cell in the 20x5 table, in place, preserving\n// each cell's run formatting (font/size) by rewriting the whole values grid.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst oldGrid = [[\"46+0=\", \"50+37=\", \"13+19=\", \"96-20=\", \"65+10=\"], [\"27+60=\", \"9+32=\", \"85-20=\", \"72+23=\", \"25+7=\"], [\"75-40=\", \"56-52=\", \"44+52=\", \"18+5=\", \"34+41=\"], [\"41+47=\", \"35+33=\", \"7+65=\", \"4+45=\", \"76-72=\"], [\"26+5=\", \"79-68=\", \"24+62=\", \"52+3=\", \"98-89=\"], [\"80-60=\", \"55-5=\", \"96-71=\", \"54-31=\", \"64-13=\"], [\"17+29=\", \"4+72=\", \"94-18=\", \"67-52=\", \"40-12=\"], [\"99-47=\", \"36-27=\", \"40+4=\", \"64-31=\", \"67+5=\"], [\"82-52=\", \"83+6=\", \"74+4=\", \"47+19=\", \"80-64=\"], [\"68+28=\", \"8+30=\", \"79-64=\", \"55-10=\", \"51-9=\"], [\"26-8=\", \"88-56=\", \"51+22=\", \"86-29=\", \"71+9=\"], [\"67+25=\", \"89-28=\", \"46+9=\", \"43+17=\", \"58-22=\"], [\"92-10=\", \"7+27=\", \"84-9=\", \"44-34=\", \"7-6=\"], [\"23-3=\", \"90-54=\", \"30+49=\", \"81-55=\", \"5+59=\"], [\"6+12=\", \"3-3=\", \"27+10=\", \"29+9=\", \"26+9=\"], [\"87-80=\", \"71-49=\", \"68+19=\", \"38-29=\", \"60+15=\"], [\"17+28=\", \"14+68=\", \"84-68=\", \"40+24=\", \"84-8=\"], [\"89-80=\", \"29-14=\", \"35-7=\", \"14+70=\", \"28+67=\"], [\"93-35=\", \"36+11=\", \"25+52=\", \"24+45=\", \"1+92=\"], [\"1+90=\", \"20-5=\", \"78+12=\", \"8+34=\", \"99-33=\"]];\nconst newGrid = [\n  [\n    \"12-9=\",\n    \"8+53=\",\n    \"87-56=\",\n    \"11+44=\",\n    \"34+15=\"\n  ],\n  [\n    \"7+11=\",\n    \"11+46=\",\n    \"79-70=\",\n    \"27-11=\",\n    \"86-48=\"\n  ],\n  [\n    \"10+50=\",\n    \"9+41=\",\n    \"24+18=\",\n    \"46-4=\",\n    \"71-5=\"\n  ],\n  [\n    \"19-16=\",\n    \"27+31=\",\n    \"38-7=\",\n    \"50+19=\",\n    \"37+24=\"\n  ],\n  [\n    \"23+3=\",\n    \"8+79=\",\n    \"68-1=\",\n    \"1+85=\",\n    \"98-45=\"\n  ],\n  [\n    \"53-16=\",\n    \"34-2=\",\n    \"16+75=\",\n    \"74-5=\",\n    \"43+44=\"\n  ],\n  [\n    \"22+28=\",\n    \"18-12=\",\n    \"35+42=\",\n    \"57-3=\",\n    \"87-85=\"\n  ],\n  [\n    \"55-15=\",\n    \"29+58=\",\n    \"51+14=\",\n    \"37-19=\",\n    \"24-3=\"\n  ],\n  [\n    \"33+66=\",\n    \"56-8=\",\n    \"8+45=\",\n    \"51+13=\",\n    \"75+0=\"\n  ],\n  [\n    \"59-24=\",\n    \"88-72=\",\n    \"78-57=\",\n    \"44+5=\",\n    \"41+24=\"\n  ],\n  [\n    \"0+71=\",\n    \"76-6=\",\n    \"12+61=\",\n    \"53-23=\",\n    \"90-21=\"\n  ],\n  [\n    \"28-24=\",\n    \"53-18=\",\n    \"81-75=\",\n    \"50+17=\",\n    \"79-59=\"\n  ],\n  [\n    \"98-23=\",\n    \"22+44=\",\n    \"33+58=\",\n    \"33+26=\",\n    \"23+72=\"\n  ],\n  [\n    \"96-54=\",\n    \"29+50=\",\n    \"8+36=\",\n    \"99-79=\",\n    \"72-52=\"\n  ],\n  [\n    \"42+46=\",\n    \"44+6=\",\n    \"67-42=\",\n    \"41-13=\",\n    \"80+0=\"\n  ],\n  [\n    \"55-47=\",\n    \"36-24=\",\n    \"6+51=\",\n    \"32+48=\",\n    \"90+2=\"\n  ],\n  [\n    \"6+64=\",\n    \"10+29=\",\n    \"31+18=\",\n    \"69-12=\",\n    \"48-47=\"\n  ],\n  [\n    \"61-17=\",\n    \"2+2=\",\n    \"90-49=\",\n    \"29+3=\",\n    \"36+9=\"\n  ],\n  [\n    \"26-12=\",\n    \"64+20=\",\n    \"7+77=\",\n    \"33-16=\",\n    \"61-11=\"\n  ],\n  [\n    \"84-1=\",\n    \"63+22=\",\n    \"99-34=\",\n    \"6+30=\",\n    \"34+48=\"\n  ]\n];\n\n// Build the replacement grid from the table's current values so the write\n// only touches cells that still hold an expected \"before\" value.\nconst currentGrid = table.values;\nconst resultGrid = currentGrid.map((row, r) => row.map((cell, c) => {\n  if (cell === oldGrid[r][c]) {\n    return newGrid[r][c];\n  }\n  return cell;\n}));\n\ntable.values = resultGrid;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$find = $d.Content.Find\n\n$replacements = @(\n    @(\"2025-03-13 Thursday\", \"2025-03-14 Friday\"),\n    @(\"46+0=\", \"12-9=\"),\n    @(\"50+37=\", \"8+53=\"),\n    @(\"13+19=\", \"87-56=\"),\n    @(\"96-20=\", \"11+44=\"),\n    @(\"65+10=\", \"34+15=\"),\n    @(\"27+60=\", \"7+11=\"),\n    @(\"9+32=\", \"11+46=\"),\n    @(\"85-20=\", \"79-70=\"),\n    @(\"72+23=\", \"27-11=\"),\n    @(\"25+7=\", \"86-48=\"),\n    @(\"75-40=\", \"10+50=\"),\n    @(\"56-52=\", \"9+41=\"),\n    @(\"44+52=\", \"24+18=\"),\n    @(\"18+5=\", \"46-4=\"),\n    @(\"34+41=\", \"71-5=\"),\n    @(\"41+47=\", \"19-16=\"),\n    @(\"35+33=\", \"27+31=\"),\n    @(\"7+65=\", \"38-7=\"),\n    @(\"4+45=\", \"50+19=\"),\n    @(\"76-72=\", \"37+24=\"),\n    @(\"26+5=\", \"23+3=\"),\n    @(\"79-68=\", \"8+79=\"),\n    @(\"24+62=\", \"68-1=\"),\n    @(\"52+3=\", \"1+85=\"),\n    @(\"98-89=\", \"98-45=\"),\n    @(\"80-60=\", \"53-16=\"),\n    @(\"55-5=\", \"34-2=\"),\n    @(\"96-71=\", \"16+75=\"),\n    @(\"54-31=\", \"74-5=\"),\n    @(\"64-13=\", \"43+44=\"),\n    @(\"17+29=\", \"22+28=\"),\n    @(\"4+72=\", \"18-12=\"),\n    @(\"94-18=\", \"35+42=\"),\n    @(\"67-52=\", \"57-3=\"),\n    @(\"40-12=\", \"87-85=\"),\n    @(\"99-47=\", \"55-15=\"),\n    @(\"36-27=\", \"29+58=\"),\n    @(\"40+4=\", \"51+14=\"),\n    @(\"64-31=\", \"37-19=\"),\n    @(\"67+5=\", \"24-3=\"),\n    @(\"82-52=\", \"33+66=\"),\n    @(\"83+6=\", \"56-8=\"),\n    @(\"74+4=\", \"8+45=\"),\n    @(\"47+19=\", \"51+13=\"),\n    @(\"80-64=\", \"75+0=\"),\n    @(\"68+28=\", \"59-24=\"),\n    @(\"8+30=\", \"88-72=\"),\n    @(\"79-64=\", \"78-57=\"),\n    @(\"55-10=\", \"44+5=\"),\n    @(\"51-9=\", \"41+24=\"),\n    @(\"26-8=\", \"0+71=\"),\n    @(\"88-56=\", \"76-6=\"),\n    @(\"51+22=\", \"12+61=\"),\n    @(\"86-29=\", \"53-23=\"),\n    @(\"71+9=\", \"90-21=\"),\n    @(\"67+25=\", \"28-24=\"),\n    @(\"89-28=\", \"53-18=\"),\n    @(\"46+9=\", \"81-75=\"),\n    @(\"43+17=\", \"50+17=\"),\n    @(\"58-22=\", \"79-59=\"),\n    @(\"92-10=\", \"98-23=\"),\n    @(\"7+27=\", \"22+44=\"),\n    @(\"84-9=\", \"33+58=\"),\n    @(\"44-34=\", \"33+26=\"),\n    @(\"7-6=\", \"23+72=\"),\n    @(\"23-3=\", \"96-54=\"),\n    @(\"90-54=\", \"29+50=\"),\n    @(\"30+49=\", \"8+36=\"),\n    @(\"81-55=\", \"99-79=\"),\n    @(\"5+59=\", \"72-52=\"),\n    @(\"6+12=\", \"42+46=\"),\n    @(\"3-3=\", \"44+6=\"),\n    @(\"27+10=\", \"67-42=\"),\n    @(\"29+9=\", \"41-13=\"),\n    @(\"26+9=\", \"80+0=\"),\n    @(\"87-80=\", \"55-47=\"),\n    @(\"71-49=\", \"36-24=\"),\n    @(\"68+19=\", \"6+51=\"),\n    @(\"38-29=\", \"32+48=\"),\n    @(\"60+15=\", \"90+2=\"),\n    @(\"17+28=\", \"6+64=\"),\n    @(\"14+68=\", \"10+29=\"),\n    @(\"84-68=\", \"31+18=\"),\n    @(\"40+24=\", \"69-12=\"),\n    @(\"84-8=\", \"48-47=\"),\n    @(\"89-80=\", \"61-17=\"),\n    @(\"29-14=\", \"2+2=\"),\n    @(\"35-7=\", \"90-49=\"),\n    @(\"14+70=\", \"29+3=\"),\n    @(\"28+67=\", \"36+9=\"),\n    @(\"93-35=\", \"26-12=\"),\n    @(\"36+11=\", \"64+20=\"),\n    @(\"25+52=\", \"7+77=\"),\n    @(\"24+45=\", \"33-16=\"),\n    @(\"1+92=\", \"61-11=\"),\n    @(\"1+90=\", \"84-1=\"),\n    @(\"20-5=\", \"63+22=\"),\n    @(\"78+12=\", \"99-34=\"),\n    @(\"8+34=\", \"6+30=\"),\n    @(\"99-33=\", \"34+48=\"),\n)\n\nforeach ($pair in $replacements) {\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $true, $false, $false, $false, $true, 1, $false, $pair[1], 2) | Out-Null\n}"}
